$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -eq $null) { continue }
    if ($val -notlike "*,*") { continue }

    $parts = $val.Split(",")
    $trimmed = @()
    foreach ($p in $parts) {
        $trimmed += $p.Trim()
    }

    $reversed = @()
    for ($i = $trimmed.Count - 1; $i -ge 0; $i--) {
        $reversed += $trimmed[$i]
    }

    $newVal = $reversed -join ", "
    $cell.Value2 = $newVal
}
